# Commit: "Changing Heat storage power to storage ratio from 20 to 10.
# For this reason, the charging and discharging capacities were updated -
# The storage capacities (MWh) remain the same."
#
# Column H on the "unitdata" sheet is capacity_output1 (the charging /
# discharging power capacity) for the "Heat storage charger" and
# "Heat storage discharger" unit rows. Halving the power-to-energy ratio
# (20 -> 10) doubles that power capacity while leaving the energy storage
# capacity (tracked elsewhere, unchanged here) as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("unitdata")

# Rows for the Heat storage charger/discharger units in both the 2025
# ("National Trends") and 2035 blocks of the table.
$rows = @(31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, `
          142, 143, 144, 145, 146, 147, 148, 149, 150, 151, 152, 153)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 8)   # Column H = capacity_output1
    $current = $cell.Value2
    $cell.Value = $current * 2
}

# Refresh the sheet's active view/selection to match the latest edits.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 143
$ws.Range("N13").Select()
